# Update "Förändrad" (column C) date values from 2023-10-08 (45207) to 2023-10-09 (45208)
# for data rows 2 through 15 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}
